$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4 values (E4:K4) - style stays the same
$ws.Range("E4").Value = 11628
$ws.Range("F4").Value = 10167
$ws.Range("G4").Value = 8066
$ws.Range("H4").Value = 8180
$ws.Range("I4").Value = 7958
$ws.Range("J4").Value = 8629
$ws.Range("K4").Value = 9072

# Update row 5 values (E5:K5) and change style (remove bottom border) to match F4 style
$ws.Range("E5").Value = 1791
$ws.Range("F5").Value = 2133
$ws.Range("G5").Value = 2031
$ws.Range("H5").Value = 2059
$ws.Range("I5").Value = 2281
$ws.Range("J5").Value = 2809
$ws.Range("K5").Value = 3592

# Copy the style from F4 (no bottom border) onto E5:K5 to match the diff's style change (s=15 -> s=7)
$ws.Range("F4").Copy()
$ws.Range("E5:K5").PasteSpecial(-4122)  # xlPasteFormats

# Update selection to A3
$ws.Range("A3").Select()
